$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 42
$ws.Range("D2").Value = "Automation3"
$ws.Range("E2").Select()
